$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Status: draft -> active
$ws.Range("B6").Value = "active"

# Experimental: true -> false
# (plain assignment of the literal "false" gets auto-coerced to a Boolean by
#  Excel; going through a text formula + paste-special-values keeps it a
#  genuine text/string cell, matching the original "true" string cell type.)
$ws.Range("B7").Formula = "=""false"""
$ws.Range("B7").Copy()
$ws.Range("B7").PasteSpecial(-4163)

# Date: 2024-12-13T10:55:58-03:00 -> 2023-01-15
# (same reasoning: keep it text, not an Excel date serial number.)
$ws.Range("B8").Formula = "=""2023-01-15"""
$ws.Range("B8").Copy()
$ws.Range("B8").PasteSpecial(-4163)
